# PowerShell Excel COM-interop script
# Applies the Betfair Back/Lay odds updates (commit: "Atualizando o arquivo XLSX")
# to the single data sheet (ActiveSheet) of the workbook.
# Each line below updates one cell's numeric value to match the new odds
# published in the source diff; only the 165 cells that actually changed
# in the diff are touched, rows/cells are left untouched otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 1.01
$ws.Range("K2").Value = 950
$ws.Range("O2").Value = 1.26
$ws.Range("Q2").Value = 1.26
$ws.Range("S2").Value = 1.26

# Row 3
$ws.Range("J3").Value = 4.2
$ws.Range("Q3").Value = 1.73
$ws.Range("S3").Value = 2.82
$ws.Range("U3").Value = 2.18

# Row 4
$ws.Range("AM4").Value = 280
$ws.Range("G4").Value = 5.9
$ws.Range("J4").Value = 3.05
$ws.Range("L4").Value = 1.55
$ws.Range("N4").Value = 2.5
$ws.Range("T4").Value = 2.18

# Row 5
$ws.Range("AF5").Value = 65
$ws.Range("AG5").Value = 980
$ws.Range("AJ5").Value = 300
$ws.Range("AK5").Value = 150
$ws.Range("AN5").Value = 220
$ws.Range("F5").Value = 6.8
$ws.Range("G5").Value = 8.199999999999999
$ws.Range("H5").Value = 1.53
$ws.Range("J5").Value = 4.1

# Row 6
$ws.Range("AD6").Value = 12.5
$ws.Range("G6").Value = 5.5
$ws.Range("S6").Value = 3.25
$ws.Range("W6").Value = 1.22
$ws.Range("Y6").Value = 10.5

# Row 8
$ws.Range("F8").Value = 1.96
$ws.Range("T8").Value = 1.69
$ws.Range("V8").Value = 1.28

# Row 9
$ws.Range("AA9").Value = 30
$ws.Range("AC9").Value = 7.8
$ws.Range("AE9").Value = 24
$ws.Range("AF9").Value = 26
$ws.Range("AG9").Value = 14.5
$ws.Range("AH9").Value = 15.5
$ws.Range("AI9").Value = 36
$ws.Range("AJ9").Value = 80
$ws.Range("AK9").Value = 40
$ws.Range("AL9").Value = 60
$ws.Range("AM9").Value = 75
$ws.Range("AN9").Value = 36
$ws.Range("AO9").Value = 16.5
$ws.Range("F9").Value = 3.45
$ws.Range("G9").Value = 3.6
$ws.Range("H9").Value = 2.3
$ws.Range("I9").Value = 2.34
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 3.6
$ws.Range("L9").Value = 1.38
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 4.2
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.42
$ws.Range("T9").Value = 1.68
$ws.Range("U9").Value = 2.32
$ws.Range("V9").Value = 1.74
$ws.Range("W9").Value = 1.39
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 11.5
$ws.Range("Z9").Value = 18.5

# Row 10
$ws.Range("G10").Value = 4.9
$ws.Range("H10").Value = 1.79
$ws.Range("S10").Value = 2.44
$ws.Range("W10").Value = 1.23

# Row 11
$ws.Range("AD11").Value = 65
$ws.Range("AE11").Value = 290
$ws.Range("AH11").Value = 36
$ws.Range("AM11").Value = 220
$ws.Range("N11").Value = 5.5
$ws.Range("O11").Value = 1.17
$ws.Range("P11").Value = 2.56
$ws.Range("Q11").Value = 1.52
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 2.28
$ws.Range("T11").Value = 2.06
$ws.Range("U11").Value = 1.75
$ws.Range("W11").Value = 4.3
$ws.Range("X11").Value = 36

# Row 12
$ws.Range("F12").Value = 2.38
$ws.Range("K12").Value = 3.75
$ws.Range("V12").Value = 1.43

# Row 13
$ws.Range("F13").Value = 1.41
$ws.Range("G13").Value = 1.66
$ws.Range("K13").Value = 6.6
$ws.Range("L13").Value = 1.4
$ws.Range("N13").Value = 2.36
$ws.Range("R13").Value = 1.15
$ws.Range("S13").Value = 3
$ws.Range("T13").Value = 1.04
$ws.Range("U13").Value = 1.04
$ws.Range("W13").Value = 2.5

# Row 14
$ws.Range("AB14").Value = 18
$ws.Range("AO14").Value = 21
$ws.Range("H14").Value = 2.78
$ws.Range("Y14").Value = 20

# Row 15
$ws.Range("H15").Value = 1.4
$ws.Range("I15").Value = 1.41
$ws.Range("J15").Value = 5.2
$ws.Range("K15").Value = 5.4
$ws.Range("V15").Value = 3.4

# Row 16
$ws.Range("L16").Value = 1.21
$ws.Range("S16").Value = 2.04
$ws.Range("T16").Value = 1.54

# Row 17
$ws.Range("AA17").Value = 46
$ws.Range("AB17").Value = 19
$ws.Range("AC17").Value = 12
$ws.Range("AD17").Value = 15
$ws.Range("AE17").Value = 30
$ws.Range("AF17").Value = 26
$ws.Range("AG17").Value = 15.5
$ws.Range("AH17").Value = 18.5
$ws.Range("AI17").Value = 38
$ws.Range("AJ17").Value = 48
$ws.Range("AK17").Value = 32
$ws.Range("AL17").Value = 980
$ws.Range("AM17").Value = 75
$ws.Range("AN17").Value = 19.5
$ws.Range("AO17").Value = 18.5
$ws.Range("G17").Value = 2.78
$ws.Range("H17").Value = 2.46
$ws.Range("I17").Value = 2.7
$ws.Range("J17").Value = 3.9
$ws.Range("N17").Value = 5.1
$ws.Range("O17").Value = 1.2
$ws.Range("Q17").Value = 1.6
$ws.Range("R17").Value = 1.56
$ws.Range("S17").Value = 2.48
$ws.Range("T17").Value = 1.55
$ws.Range("U17").Value = 2.48
$ws.Range("V17").Value = 1.59
$ws.Range("W17").Value = 1.56
$ws.Range("X17").Value = 29
$ws.Range("Y17").Value = 18.5
$ws.Range("Z17").Value = 25

# Row 18
$ws.Range("G18").Value = 1.86
$ws.Range("J18").Value = 3.45
$ws.Range("Q18").Value = 2.22
$ws.Range("S18").Value = 3.8
$ws.Range("U18").Value = 1.75
$ws.Range("V18").Value = 1.14
$ws.Range("W18").Value = 2.16

# Row 19
$ws.Range("V19").Value = 1.46
$ws.Range("X19").Value = 12

# Row 20
$ws.Range("F20").Value = 1.85
$ws.Range("G20").Value = 1.86
$ws.Range("H20").Value = 5.5
$ws.Range("I20").Value = 5.6
$ws.Range("V20").Value = 1.21
$ws.Range("W20").Value = 2.16

# Row 21
$ws.Range("AA21").Value = 540
$ws.Range("AD21").Value = 44
$ws.Range("AE21").Value = 220
$ws.Range("AM21").Value = 210
$ws.Range("AO21").Value = 290
$ws.Range("F21").Value = 1.33
$ws.Range("G21").Value = 1.35
$ws.Range("I21").Value = 12.5
$ws.Range("Q21").Value = 1.77
$ws.Range("V21").Value = 1.08
$ws.Range("W21").Value = 3.85
$ws.Range("Y21").Value = 38
$ws.Range("Z21").Value = 120
